# RCM-exp.xlsx -- add the "Bandwidth after RCM ordering" sheets.
#
# Sheet1 keeps its data, but row 11's labels get centred (new style) and the
# selection moves from E18 to C6 (and Sheet1 is no longer the active tab).
# Three new sheets are appended: "new run", "Sheet3", "Sheet4" -- each a
# variant of the Sheet1 layout with extra RCM-ordering bandwidth numbers,
# ending with Sheet4 selected.

$wb = $excel.ActiveWorkbook
$s1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1 tweaks
# ---------------------------------------------------------------------
$s1.Range("B11,C11,D11,G11,H11,I11,J11").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Add the three new worksheets, in order, right after Sheet1
# ---------------------------------------------------------------------
$s2 = $wb.Worksheets.Add($null, $s1)
$s2.Name = "new run"

$s3 = $wb.Worksheets.Add($null, $s2)
$s3.Name = "Sheet3"

$s4 = $wb.Worksheets.Add($null, $s3)
$s4.Name = "Sheet4"

# ---------------------------------------------------------------------
# Column widths -- match Sheet1's column layout (A / B:C / D:F / G / H)
# ---------------------------------------------------------------------
foreach ($ws in @($s2, $s3)) {
    $ws.Range("A1").ColumnWidth = 16.833333333333332
    $ws.Range("B1:C1").ColumnWidth = 13.666666666666666
    $ws.Range("D1:F1").ColumnWidth = 13.333333333333334
    $ws.Range("G1").ColumnWidth = 16.5
    $ws.Range("H1").ColumnWidth = 13.166666666666666
}

# ---------------------------------------------------------------------
# Page setup -- match Sheet1's margins; sheets 2 & 3 also get portrait
# orientation explicitly set (Sheet4 is left at engine defaults)
# ---------------------------------------------------------------------
foreach ($ws in @($s2, $s3, $s4)) {
    $ps = $ws.PageSetup
    $ps.LeftMargin = $excel.InchesToPoints(0.75)
    $ps.RightMargin = $excel.InchesToPoints(0.75)
    $ps.TopMargin = $excel.InchesToPoints(1)
    $ps.BottomMargin = $excel.InchesToPoints(1)
    $ps.HeaderMargin = $excel.InchesToPoints(0.5)
    $ps.FooterMargin = $excel.InchesToPoints(0.5)
}
$s2.PageSetup.Orientation = 1
$s3.PageSetup.Orientation = 1

# =======================================================================
# "new run" sheet (sheet2)
# =======================================================================
$s2.Range("B1").Value = "hugetric-00020"
$s2.Range("D1").Value = "dielFilterV3real"
$s2.Range("G1").Value = "delaunay_n24.mtx"
$s2.Range("K1").Value = "hugetric-00020"
$s2.Range("L1").Value = "dielFilterV3real"
$s2.Range("M1").Value = "delaunay_n24.mtx"

$s2.Range("A2").Value = "HSL time"
$s2.Range("B2").Value = 3.83
$s2.Range("D2").Value = 0.74
$s2.Range("G2").Value = 4.1
$s2.Range("J2").Value = "HSL time"
$s2.Range("K2").Value = 3.83
$s2.Range("L2").Value = 0.74
$s2.Range("M2").Value = 4.1

$s2.Range("A3").Value = "pseudo-diameter"
$s2.Range("B3").Value = 3661
$s2.Range("D3").Value = 84
$s2.Range("G3").Value = 1720
$s2.Range("J3").Value = "pseudo-diameter"
$s2.Range("K3").Value = 3661
$s2.Range("L3").Value = 84
$s2.Range("M3").Value = 1720

$s2.Range("A5").Value = "cores"
$s2.Range("J5").Value = "cores"

$s2.Range("A6").Value = 1
$s2.Range("D6").Value = 27.845
$s2.Range("E6").Value = 9.06414
$s2.Range("F6").Value = 0.694179

$s2.Range("A7").Value = 4
$s2.Range("D7").Value = 4.86092
$s2.Range("E7").Value = 2.36754
$s2.Range("F7").Value = 0.258353

$s2.Range("A8").Value = 16
$s2.Range("D8").Value = 2.27935
$s2.Range("E8").Value = 0.746885
$s2.Range("F8").Value = 0.134313
$s2.Range("J8").Value = 24
$s2.Range("L8").Value = 2.36128
$s2.Range("M8").Value = 80.4994

$s2.Range("A9").Value = 64
$s2.Range("D9").Value = 0.90816
$s2.Range("E9").Value = 0.25439
$s2.Range("F9").Value = 0.747488
$s2.Range("J9").Value = 96
$s2.Range("L9").Value = 0.852554
$s2.Range("M9").Value = 25.3132

$s2.Range("A10").Value = 256
$s2.Range("D10").Value = 0.416558
$s2.Range("E10").Value = 0.1034
$s2.Range("F10").Value = 9.37438
$s2.Range("J10").Value = 384
$s2.Range("L10").Value = 0.341284
$s2.Range("M10").Value = 13.7091

$s2.Range("A11").Value = 1024
$s2.Range("D11").Value = 0.460892
$s2.Range("J11").Value = 1536
$s2.Range("L11").Value = 0.148418
$s2.Range("M11").Value = 6.24828

$s2.Range("B13").Value = "unthreaded"
$s2.Range("J13").Value = "6 threads per MPI process"

# centred header/footer cells
$s2.Range("B1,C1,D1,E1,F1,G1,H1").HorizontalAlignment = -4108
$s2.Range("B13,C13,D13,E13,F13,G13,J13,K13,L13,M13").HorizontalAlignment = -4108

# merges (insertion order matters for the XML's mergeCells listing)
$s2.Range("B13:G13").Merge()
$s2.Range("J13:M13").Merge()
$s2.Range("D1:E1").Merge()
$s2.Range("G1:H1").Merge()
$s2.Range("B1:C1").Merge()

$s2.Range("F12").Select()

# =======================================================================
# "Sheet3"
# =======================================================================
$s3.Range("B1").Value = "hugetric-00020"
$s3.Range("D1").Value = "dielFilterV3real"
$s3.Range("G1").Value = "delaunay_n24.mtx"

$s3.Range("A5").Value = "cores"

$s3.Range("A6").Value = 1
$s3.Range("A7").Value = 4
$s3.Range("A8").Value = 16
$s3.Range("A9").Value = 64

$s3.Range("A10").Value = 256
$s3.Range("D10").Value = 0.519198
$s3.Range("E10").Value = 0.172929
$s3.Range("F10").Value = 0.17516
$s3.Range("G10").Value = 25.7102
$s3.Range("H10").Value = 6.48163
$s3.Range("I10").Value = 4.27308

$s3.Range("A11").Value = 1024
$s3.Range("D11").Value = 0.349325
$s3.Range("E11").Value = 0.0579925
$s3.Range("F11").Value = 1.30855
$s3.Range("G11").Value = 11.2601
$s3.Range("H11").Value = 3.40651
$s3.Range("I11").Value = 58.6354

$s3.Range("B13").Value = "unthreaded"

$s3.Range("B1,C1,D1,E1,F1,G1,H1").HorizontalAlignment = -4108
$s3.Range("B13,C13,D13,E13,F13,G13,J13,K13,L13,M13").HorizontalAlignment = -4108

$s3.Range("B1:C1").Merge()
$s3.Range("D1:E1").Merge()
$s3.Range("G1:H1").Merge()
$s3.Range("B13:G13").Merge()
$s3.Range("J13:M13").Merge()

$s3.Range("D10").Select()

# =======================================================================
# "Sheet4"
# =======================================================================
$s4.Range("E3").Value = 1
$s4.Range("F3").Value = 90
$s4.Range("H3").Value = 30

$s4.Range("E4").Value = 16
$s4.Range("F4").Value = 17.945
$s4.Range("H4").Value = 8

$s4.Range("E5").Value = 256
$s4.Range("F5").Value = 5

$s4.Range("L19").Select()

# Sheet4 ends up the active tab/selected sheet.
$s4.Select()

# Finally, move the Sheet1 selection (it was E18, now C6) -- done last so
# it doesn't get clobbered by the sheet-activation calls above.
$s1.Range("C6").Select()
$s4.Select()
